$d = $word.ActiveDocument
$wns = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

# ---------------------------------------------------------------------------
# 1) "Sportify" + " (JavaScript)" -> single run "Sportify (JavaScript)"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Sportify (JavaScript)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraSportify = $r.Paragraphs(1)
$xmlSportify = "<w:p $wns><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"7`"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:t>Sportify (JavaScript)</w:t></w:r></w:p>"
$paraSportify.Range.InsertXML($xmlSportify)

# ---------------------------------------------------------------------------
# 2) Empty paragraph right before the "Maquette-du-site-Voiture-de-luxe" link
#    gets turned into a numbered list item: Site Web « Voiture de luxe »
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("HabitationPlus.png.png", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraHab = $r.Paragraphs(1)
$spacer = $paraHab.Next()
$target1 = $spacer.Next()
# Prime a fresh list definition then overwrite with the real content / numId.
$target1.Range.ListFormat.ApplyNumberDefault()
$xml1 = "<w:p $wns><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"9`"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:t>Site Web</w:t></w:r><w:r><w:t xml:space=`"preserve`"> &#171; </w:t></w:r><w:r><w:t>Voiture de luxe</w:t></w:r><w:r><w:t> &#187;</w:t></w:r></w:p>"
$target1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 3) Insert a brand-new list item right before the empty paragraph that
#    follows "Voiture de luxe.png": Site Web « Voiture de luxe » sous format mobile
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Voiture de luxe.png", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraVoiture = $r.Paragraphs(1)
$target2 = $paraVoiture.Next()
$insertionRange = $target2.Range
$insertionRange.Collapse(1)
$insertionRange.InsertParagraphBefore()

$r = $d.Content
$r.Find.Execute("Voiture de luxe.png", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraVoiture = $r.Paragraphs(1)
$newPara = $paraVoiture.Next()
$newPara.Range.ListFormat.ApplyNumberDefault()
$xml2 = "<w:p $wns><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"9`"/></w:numPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:t xml:space=`"preserve`">Site Web </w:t></w:r><w:r><w:t>&#171; </w:t></w:r><w:r><w:t>Voiture de luxe</w:t></w:r><w:r><w:t> &#187;</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>sous format mobile</w:t></w:r></w:p>"
$newPara.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 4) "Voiture de luxe(mobile)" + " " -> single run "Voiture de luxe(mobile) "
#    (the "Image : " and ".png" runs stay as they are)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Voiture de luxe(mobile)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraMobile = $r.Paragraphs(1)
$xml3 = "<w:p $wns><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:t xml:space=`"preserve`">Image : </w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space=`"preserve`">Voiture de luxe(mobile) </w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>.png</w:t></w:r></w:p>"
$paraMobile.Range.InsertXML($xml3)
